# Adapt column header formatting to respective input file names (#7)
# - rename "<header>_old"  -> "<header>_FV2310"
# - rename "<header>_new"  -> "<header>_FV2404"
# - turn the data range into an Excel Table ("Table1")
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldSuffix = "_old"
$newSuffix = "_new"
$fv2310 = "_FV2310"
$fv2404 = "_FV2404"

# Rename the header cells (row 1) in place, preserving order/position.
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $text = $cell.Value()

    if ($text -like "*$oldSuffix") {
        $base = $text.Substring(0, $text.Length - $oldSuffix.Length)
        $cell.Value = "$base$fv2310"
    }
    elseif ($text -like "*$newSuffix") {
        $base = $text.Substring(0, $text.Length - $newSuffix.Length)
        $cell.Value = "$base$fv2404"
    }
}

# Turn the used range into a native Excel table ("Table1") with the
# (freshly renamed) header row as column headers.
$dataRange = $ws.Range("A1:U72")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"
$table.TableStyle = $null

# Freeze the header row (pane split below row 1).
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
